$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray empty "Proxy" cell at F27 (column F has no data for this row)
$ws.Range("F27").ClearContents()

# Append new account rows (28-42)
# Row 28
$ws.Range("A28").Value = 'Emma Williams'
$ws.Range("B28").Value = 'yotixog229@bauscn.com'
$ws.Range("C28").Value = 'Pass1231'
$ws.Range("D28").Value = '28 April 1982'
$ws.Range("E28").Value = 'female'
$ws.Range("G28").Value = '2025-04-19 21:47:28'

# Row 29
$ws.Range("A29").Value = 'Emma Williams'
$ws.Range("B29").Value = 'yotixog229@bauscn.com'
$ws.Range("C29").Value = 'Pass1231'
$ws.Range("D29").Value = '28 April 1982'
$ws.Range("E29").Value = 'female'
$ws.Range("G29").Value = '2025-04-19 21:47:28'

# Row 30
$ws.Range("A30").Value = 'John Johnson'
$ws.Range("B30").Value = 'woroko1536@f5url.com'
$ws.Range("C30").Value = 'Pass9639'
$ws.Range("D30").Value = '21 October 1984'
$ws.Range("E30").Value = 'female'
$ws.Range("G30").Value = '2025-04-19 21:59:02'

# Row 31
$ws.Range("A31").Value = 'John Johnson'
$ws.Range("B31").Value = 'woroko1536@f5url.com'
$ws.Range("C31").Value = 'Pass9639'
$ws.Range("D31").Value = '21 October 1984'
$ws.Range("E31").Value = 'female'
$ws.Range("G31").Value = '2025-04-19 21:59:02'

# Row 32
$ws.Range("A32").Value = 'John Johnson'
$ws.Range("B32").Value = 'poyef47285@cotigz.com'
$ws.Range("C32").Value = 'Pass9639'
$ws.Range("D32").Value = '21 October 1984'
$ws.Range("E32").Value = 'female'
$ws.Range("G32").Value = '2025-04-19 21:59:45'

# Row 33
$ws.Range("A33").Value = 'John Johnson'
$ws.Range("B33").Value = 'poyef47285@cotigz.com'
$ws.Range("C33").Value = 'Pass9639'
$ws.Range("D33").Value = '21 October 1984'
$ws.Range("E33").Value = 'female'
$ws.Range("G33").Value = '2025-04-19 21:59:45'

# Row 34
$ws.Range("A34").Value = 'James Williams'
$ws.Range("B34").Value = 'james.jones47@yahoo.com'
$ws.Range("C34").Value = 'Pass3195'
$ws.Range("D34").Value = '12 June 1989'
$ws.Range("E34").Value = 'male'
$ws.Range("G34").Value = '2025-04-19 22:02:00'

# Row 35
$ws.Range("A35").Value = 'Olivia Johnson'
$ws.Range("B35").Value = 'nigip97395@linxues.com'
$ws.Range("C35").Value = 'Pass5396'
$ws.Range("D35").Value = '17 September 1981'
$ws.Range("E35").Value = 'female'
$ws.Range("G35").Value = '2025-04-19 22:03:52'

# Row 36
$ws.Range("A36").Value = 'Olivia Johnson'
$ws.Range("B36").Value = 'nigip97395@linxues.com'
$ws.Range("C36").Value = 'Pass5396'
$ws.Range("D36").Value = '17 September 1981'
$ws.Range("E36").Value = 'female'
$ws.Range("G36").Value = '2025-04-19 22:03:52'

# Row 37
$ws.Range("A37").Value = 'John Brown'
$ws.Range("B37").Value = 'mosakar418@f5url.com'
$ws.Range("C37").Value = 'Pass4832'
$ws.Range("D37").Value = '19 February 1989'
$ws.Range("E37").Value = 'male'
$ws.Range("G37").Value = '2025-04-19 23:32:20'

# Row 38
$ws.Range("A38").Value = 'John Brown'
$ws.Range("B38").Value = 'mosakar418@f5url.com'
$ws.Range("C38").Value = 'Pass4832'
$ws.Range("D38").Value = '19 February 1989'
$ws.Range("E38").Value = 'male'
$ws.Range("G38").Value = '2025-04-19 23:32:20'

# Row 39
$ws.Range("A39").Value = 'James Smith'
$ws.Range("B39").Value = 'jawidam390@agiuse.com'
$ws.Range("C39").Value = 'Pass1348'
$ws.Range("D39").Value = '7 September 1999'
$ws.Range("E39").Value = 'male'
$ws.Range("G39").Value = '2025-04-19 23:36:34'

# Row 40
$ws.Range("A40").Value = 'James Smith'
$ws.Range("B40").Value = 'jawidam390@agiuse.com'
$ws.Range("C40").Value = 'Pass1348'
$ws.Range("D40").Value = '7 September 1999'
$ws.Range("E40").Value = 'male'
$ws.Range("G40").Value = '2025-04-19 23:36:38'

# Row 41
$ws.Range("A41").Value = 'James Smith'
$ws.Range("B41").Value = 'wimatog888@cxnlab.com'
$ws.Range("C41").Value = 'Pass1348'
$ws.Range("D41").Value = '7 September 1999'
$ws.Range("E41").Value = 'male'
$ws.Range("G41").Value = '2025-04-19 23:37:16'

# Row 42
$ws.Range("A42").Value = 'James Smith'
$ws.Range("B42").Value = 'wimatog888@cxnlab.com'
$ws.Range("C42").Value = 'Pass1348'
$ws.Range("D42").Value = '7 September 1999'
$ws.Range("E42").Value = 'male'
$ws.Range("G42").Value = '2025-04-19 23:37:16'
